$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "Nomor HP" (column C) to hold the new
# "Nama Pegawai" field. This shifts Nomor HP, Jabatan, Nomor Kendaraan,
# Nama Kendaraan and Tipe one column to the right.
$ws.Range("C:C").Insert()

# New column header / shared string "Nama Pegawai"
$ws.Range("C2").Value = "Nama Pegawai"

# Narrow the new "Nama Biro" column slightly.
$ws.Range("B:B").ColumnWidth = 14.7

# Match the resulting alignment tweaks on the data rows by copying the
# existing right-aligned/left-aligned formats from neighbouring cells
# (keeps the style table free of duplicate/unused entries).
$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("A4").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)

$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Restore the active selection recorded in the saved workbook.
$ws.Range("D11").Select()
